$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("130, 780", 1, "130", "3631"),
    @("780", 1, "780", "3666"),
    @("130, 1073", 2, "130, 130", "5131, 4415"),
    @("780, 1073, 1105", 1, "1105", "4994"),
    @("130, 455, 780", 2, "130, 130", "5269, 6424"),
    @("455, 1073, 1105", 1, "1105", "5399"),
    @("423, 748, 780, 1073", 2, "780, 780", "5677, 5887"),
    @("423, 1073, 1105", 2, "1105, 1105", "5331, 5433"),
    @("98, 130, 455, 748, 1073", 1, "130", "6016"),
    @("98, 130, 748, 1073", 1, "130", "5582"),
    @("130, 423, 748, 1073, SF", 1, "130", "6202"),
    @("98, 130, 780, 1073", 1, "130", "6488"),
    @("130, 423, 780, 1073", 1, "130", "6561"),
    @("130, 780, 780, 1073", 1, "780", "6359"),
    @("98, 98, 130, 455, 780", 1, "130", "6727"),
    @("98, 130, 423, 1073", 1, "130", "6670"),
    @("98, 98, 455, 780", 1, "780", "6748")
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).NumberFormat = "@"
    $ws.Cells.Item($rowIndex, 3).NumberFormat = "@"
    $ws.Cells.Item($rowIndex, 4).NumberFormat = "@"

    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]

    $ws.Cells.Item($rowIndex, 1).ClearFormats()
    $ws.Cells.Item($rowIndex, 3).ClearFormats()
    $ws.Cells.Item($rowIndex, 4).ClearFormats()

    $rowIndex++
}
